$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82 (pushes existing rows 82..136 down to 83..137)
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new weekly data point
$ws.Cells.Item(82, 1).Value = 10
$ws.Cells.Item(82, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(82, 3).Value = "La Araucanía"
$ws.Cells.Item(82, 4).Value = 45236
$ws.Cells.Item(82, 5).Value = 9
$ws.Cells.Item(82, 6).Value = 300000001
$ws.Cells.Item(82, 7).Value = "Rabanito"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 70
$ws.Cells.Item(82, 11).Value = 8000
$ws.Cells.Item(82, 12).Value = 9000
$ws.Cells.Item(82, 13).Value = 8714
$ws.Cells.Item(82, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(82, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(82, 16).Value = 726
$ws.Cells.Item(82, 17).Value = 12
$ws.Cells.Item(82, 18).Value = "Hortaliza"
